# Restored from revision #2b64b5eb8e5394a6072ba2303b3be6357f67aba8.TEST
# Author: admin. Type: SAVE.
#
# Semantic change: cell C10 on the active sheet (the "R30" rule row,
# condition C2 = "hour <= max") changes its numeric value from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
